# Unit_Status.xlsx update:
#  - Added 12 new units (2I,2J,2K,2L, 3I,3J,3K,3L, 4I,4J,4K,4L)
#  - Re-sorted data rows by floor/unit (1A-1D, 2A-2L, 3A-3L, 4A-4L)
#  - Fixed several "Null"/text-date cells so Occupancy Status + Lease Date
#    display correctly (UI notification bug), and corrected some lease dates
#  - 1A's Lease Date / Lease Period now show "Vacant" (instead of being blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash the 3 existing cell styles we need to reuse into scratch cells
#     well outside the data block, BEFORE the data block gets cleared/
#     overwritten. This lets every row below pick up the same cellXfs
#     entries instead of Excel minting a brand-new style per write.
#       H1 -> numFmt "m/d/yyyy" date, right aligned   (style 1, like old D4)
#       H2 -> numFmt "m/d/yyyy" date, general aligned  (style 2, like old D3)
#       H3 -> plain/default style, for text cells      (style 0, like old C2)
$ws.Range("D4").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$tplRight   = "H1"
$tplGeneral = "H2"
$tplPlain   = "H3"

# Wipe the existing data block (row 1 header untouched) and any stale rows
# below the old range, then rewrite every row in the new order.
$ws.Range("A2:E41").Clear()

$rows = @(
     ,@(2,  "1A", 56.4,  "Vacant",   "Vacant", "Vacant", "0")
     ,@(3,  "1B", 57.58, "Occupied", 45345,    1,        "2")
     ,@(4,  "1C", 57.58, "Occupied", 45391,    1,        "1")
     ,@(5,  "1D", 58.75, "Occupied", 45429,    1,        "2")
     ,@(6,  "2A", 22.82, "Occupied", 45430,    1,        "2")
     ,@(7,  "2B", 35,    "Occupied", 45128,    1,        "1")
     ,@(8,  "2C", 22.82, "Occupied", 45444,    1,        "1")
     ,@(9,  "2D", 35,    "Occupied", 45118,    1,        "1")
     ,@(10, "2E", 22.82, "Vacant",   "Null",   "Null",   "0")
     ,@(11, "2F", 35,    "Occupied", 45023,    1,        "1")
     ,@(12, "2G", 22.82, "Occupied", 45397,    1,        "1")
     ,@(13, "2H", 35,    "Occupied", 45177,    1,        "1")
     ,@(14, "2I", 22.82, "Occupied", 45372,    1,        "1")
     ,@(15, "2J", 35,    "Occupied", 45304,    1,        "1")
     ,@(16, "2K", 22.82, "Occupied", 45299,    1,        "1")
     ,@(17, "2L", 35,    "Occupied", 45348,    1,        "1")
     ,@(18, "3A", 22.82, "Occupied", 45298,    1,        "1")
     ,@(19, "3B", 35,    "Occupied", 45148,    1,        "2")
     ,@(20, "3C", 22.82, "Occupied", 45171,    1,        "2")
     ,@(21, "3D", 35,    "Occupied", 45288,    1,        "1")
     ,@(22, "3E", 22.82, "Occupied", 45245,    1,        "1")
     ,@(23, "3F", 35,    "Occupied", 45224,    1,        "2")
     ,@(24, "3G", 22.82, "Occupied", 45137,    1,        "2")
     ,@(25, "3H", 35,    "Occupied", 45156,    1,        "2")
     ,@(26, "3I", 22.82, "Occupied", 45143,    1,        "2")
     ,@(27, "3J", 35,    "Occupied", 45219,    1,        "2")
     ,@(28, "3K", 22.82, "Occupied", 45210,    1,        "2")
     ,@(29, "3L", 35,    "Occupied", 45388,    1,        "2")
     ,@(30, "4A", 22.82, "Occupied", 45275,    1,        "1")
     ,@(31, "4B", 35,    "Occupied", 45225,    1,        "1")
     ,@(32, "4C", 22.82, "Occupied", 45237,    1,        "1")
     ,@(33, "4D", 35,    "Occupied", 45166,    1,        "1")
     ,@(34, "4E", 22.82, "Occupied", 45303,    1,        "2")
     ,@(35, "4F", 35,    "Occupied", 45275,    1,        "2")
     ,@(36, "4G", 22.82, "Occupied", 45392,    1,        "1")
     ,@(37, "4H", 35,    "Occupied", 45338,    1,        "1")
     ,@(38, "4I", 22.82, "Occupied", 45170,    1,        "1")
     ,@(39, "4J", 35,    "Occupied", 45182,    1,        "1")
     ,@(40, "4K", 22.82, "Occupied", 45320,    1,        "1")
     ,@(41, "4L", 35,    "Occupied", 45411,    1,        "1")
)

foreach ($row in $rows) {
    $r      = $row[0]
    $unit   = $row[1]
    $area   = $row[2]
    $occ    = $row[3]
    $dateV  = $row[4]
    $period = $row[5]
    $style  = $row[6]

    $ws.Cells.Item($r, 1).Value2 = $unit
    $ws.Cells.Item($r, 2).Value2 = $area
    $ws.Cells.Item($r, 3).Value2 = $occ

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    if ($style -eq "0") {
        # Text cell (Vacant / Null) - plain default style, no number format.
        $ws.Range($tplPlain).Copy()
        $dCell.PasteSpecial(-4122)
        $ws.Range($tplPlain).Copy()
        $eCell.PasteSpecial(-4122)

        $dCell.Value2 = $dateV
        $eCell.Value2 = $period
    } else {
        if ($style -eq "1") {
            $ws.Range($tplRight).Copy()
        } else {
            $ws.Range($tplGeneral).Copy()
        }
        $dCell.PasteSpecial(-4122)
        $dCell.Value2 = $dateV
        $eCell.Value2 = $period
    }
}

$excel.CutCopyMode = $false

# Remove the scratch template cells now that every row has its style.
$ws.Range("H1:H3").Clear()

# Restore the sheet-view bits captured in the source edit (zoom + selection).
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("D25").Select()
